$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: add the new Grade-7 rows using the existing column layout
#     (A = Grade, B = Subjects, C = Con) before the new "Lesson" column
#     is inserted. This matches the order in which Excel would intern
#     new shared strings as the rows were typed in.
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "English"
$ws.Range("C5").Value = "Quiz,Worksheet,Flashcards"

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Math"
$ws.Range("C6").Value = "Worksheet"

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Science"
$ws.Range("C7").Value = "Quiz,Worksheet,Flashcards"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Social"
$ws.Range("C8").Value = "Quiz,Worksheet,Flashcards"

# --- Step 2: insert a new column C ("Lesson") in front of the old "Con"
#     column, which shifts Con (and everything in it) to column D.
$ws.Range("C1").EntireColumn.Insert()

# --- Step 3: fill in the new "Lesson" column, top to bottom, for every row.
$ws.Range("C1").Value = "Lesson"
$ws.Range("C2").Value = "Lesson 1"
$ws.Range("C3").Value = "Addition"
$ws.Range("C4").Value = "Plants"
$ws.Range("C5").Value = "Tenses"
$ws.Range("C6").Value = "Ratio & Propotion"
$ws.Range("C7").Value = "Forest Our Life Line"
$ws.Range("C8").Value = "Market"

# Update the selected cell to mimic the final saved view state.
$ws.Range("D12").Select()
